# Applies the "cfast vv: add some annotation to cfast verification spreadsheet"
# edit to Sheet1 of the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Clear the scratch values in D2:E2 (row becomes empty and disappears).
$ws.Range("D2:E2").ClearContents()

# 2. D4 used to be a formula "=D2*E2"; now that D2/E2 are gone it just holds
#    the literal value it used to compute.
$ws.Range("D4").Value = 6000000

# 3. Rename "V" (A5) -> "Volume"
$ws.Range("A5").Value = "Volume"

# 4. Rename "calc pressure" (C8) -> "calculated pressure"
$ws.Range("C8").Value = "calculated pressure"

# 5. Rename "T" (A9) -> "Time"
$ws.Range("A9").Value = "Time"

# 6. Add "temperature" label in F9 (was blank)
$ws.Range("F9").Value = "temperature"

# 7. New italic note row 7: F7
$ws.Range("F7").Value = "cfast temperature and pressur columns are copied from a CFAST run"
$ws.Range("F7").Font.Italic = $true

# 8. New annotation block in columns N:Q, rows 8-13 (italic style, no text except
#    N8 and N10:N13)
$ws.Range("N8").Value = "Formulas (assuming constant fire)"
$ws.Range("N10").Value = "DP=(gamma-1)*qtotal*Time/Volume"
$ws.Range("N11").Value = "M=M0+mfire*Time"
$ws.Range("N12").Value = "E=E0+qconvec*Time"
$ws.Range("N13").Value = "T=E/(cv*M)-273.3"

$ws.Range("N8:Q13").Font.Italic = $true

# 9. Update the selected cell shown when the workbook is opened.
$ws.Range("D5").Select() | Out-Null
